# Update Leve profit-tracking numbers (currentAveragePrice*, LevePrice*, LeveProfit*)
# across all job sheets, reflecting refreshed market-board pricing data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18: You Grow, Girl
$ws.Range("H18").Value = 13495.208
$ws.Range("I18").Value = 8909.091
$ws.Range("J18").Value = 17375.77
$ws.Range("K18").Value = 8909.091
$ws.Range("L18").Value = 17375.77
$ws.Range("M18").Value = -8625.091
$ws.Range("N18").Value = -17943.77

# Row 33: Glazed and Confused
$ws.Range("H33").Value = 99.333336
$ws.Range("I33").Value = 99.333336
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 99.333336
$ws.Range("L33").ClearContents()
$ws.Range("M33").Value = 129.666664
$ws.Range("N33").Value = 0

# Row 98: The Dotted Line
$ws.Range("H98").Value = 3541
$ws.Range("I98").Value = 4571
$ws.Range("K98").Value = 4571
$ws.Range("M98").Value = -3073

# Row 101: Edge of the Arcane
$ws.Range("H101").Value = 819.1667
$ws.Range("I101").Value = 272.5
$ws.Range("J101").Value = 1092.5
$ws.Range("K101").Value = 817.5
$ws.Range("L101").Value = 3277.5
$ws.Range("M101").Value = 804.5
$ws.Range("N101").Value = -6521.5

# Row 122: Wishful Inking
$ws.Range("H122").Value = 3541
$ws.Range("I122").Value = 4571
$ws.Range("K122").Value = 13713
$ws.Range("M122").Value = -11263

# Row 125: Body over Mind
$ws.Range("H125").Value = 316.5
$ws.Range("I125").Value = 380
$ws.Range("J125").Value = 253
$ws.Range("K125").Value = 3420
$ws.Range("L125").Value = 2277
$ws.Range("M125").Value = -960
$ws.Range("N125").Value = -7197

# Row 141: Remedy for Reason
$ws.Range("H141").Value = 1079116.2
$ws.Range("I141").Value = 1335258.2
$ws.Range("J141").Value = 3320
$ws.Range("K141").Value = 4005774.6
$ws.Range("L141").Value = 9960
$ws.Range("M141").Value = -4000594.6
$ws.Range("N141").Value = -20320

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 3917.6316
$ws.Range("I32").Value = 3378.1592
$ws.Range("J32").Value = 5743.5386
$ws.Range("K32").Value = 3378.1592
$ws.Range("L32").Value = 5743.5386
$ws.Range("M32").Value = -3091.1592
$ws.Range("N32").Value = -6317.5386

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 2635.5908
$ws.Range("I61").Value = 1331.2
$ws.Range("K61").Value = 1331.2
$ws.Range("M61").Value = -1119.2

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 939.4
$ws.Range("I74").Value = 543.6957
$ws.Range("J74").Value = 5490
$ws.Range("K74").Value = 543.6957
$ws.Range("L74").Value = 5490
$ws.Range("M74").Value = 330.3043
$ws.Range("N74").Value = -7238

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 939.4
$ws.Range("I77").Value = 543.6957
$ws.Range("J77").Value = 5490
$ws.Range("K77").Value = 2718.4785
$ws.Range("L77").Value = 27450
$ws.Range("M77").Value = 1649.5215
$ws.Range("N77").Value = -36186

# Row 109: A Head of Demand
$ws.Range("H109").Value = 67290.336
$ws.Range("J109").Value = 67290.336
$ws.Range("L109").Value = 67290.336
$ws.Range("N109").Value = -70064.336

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 2635.5908
$ws.Range("I136").Value = 1331.2
$ws.Range("K136").Value = 3993.6
$ws.Range("M136").Value = -1443.6

$ws = $wb.Worksheets.Item("BSM")
# Row 132: Always Be Prepaired
$ws.Range("H132").Value = 122780
$ws.Range("J132").Value = 122780
$ws.Range("L132").Value = 122780
$ws.Range("N132").Value = -132900

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 2917.5588
$ws.Range("I134").Value = 2975.6667
$ws.Range("K134").Value = 8927.000100000001
$ws.Range("M134").Value = -6392.000100000001

$ws = $wb.Worksheets.Item("CRP")
# Row 4: A Clogful of Camaraderie
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = 0

# Row 31: Wall Not Found
$ws.Range("H31").Value = 1648.1904
$ws.Range("I31").Value = 989.1667
$ws.Range("J31").Value = 2526.889
$ws.Range("K31").Value = 989.1667
$ws.Range("L31").Value = 2526.889
$ws.Range("M31").Value = -694.1667
$ws.Range("N31").Value = -3116.889

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 1648.1904
$ws.Range("I34").Value = 989.1667
$ws.Range("J34").Value = 2526.889
$ws.Range("K34").Value = 989.1667
$ws.Range("L34").Value = 2526.889
$ws.Range("M34").Value = -787.1667
$ws.Range("N34").Value = -2930.889

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 5436193
$ws.Range("J58").Value = 2187.5
$ws.Range("L58").Value = 2187.5
$ws.Range("N58").Value = -2593.5

# Row 99: O Pine
$ws.Range("H99").Value = 1252489
$ws.Range("I99").Value = 2502224.5
$ws.Range("K99").Value = 2502224.5
$ws.Range("M99").Value = -2500726.5

# Row 126: A Better Conductor
$ws.Range("H126").Value = 1252489
$ws.Range("I126").Value = 2502224.5
$ws.Range("K126").Value = 7506673.5
$ws.Range("M126").Value = -7504203.5

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 2134.476
$ws.Range("I132").Value = 1636
$ws.Range("J132").Value = 4253
$ws.Range("K132").Value = 4908
$ws.Range("L132").Value = 12759
$ws.Range("M132").Value = -2378
$ws.Range("N132").Value = -17819

# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 1687.8214
$ws.Range("I134").Value = 1018.5909
$ws.Range("J134").Value = 4141.6665
$ws.Range("K134").Value = 3055.7727
$ws.Range("L134").Value = 12424.9995
$ws.Range("M134").Value = -520.7727
$ws.Range("N134").Value = -17494.9995

# Row 136: Turali Quality
$ws.Range("H136").Value = 5436193
$ws.Range("J136").Value = 2187.5
$ws.Range("L136").Value = 6562.5
$ws.Range("N136").Value = -11662.5

# Row 141: No Greater Treasure
$ws.Range("H141").Value = 70106.664
$ws.Range("J141").Value = 70106.664
$ws.Range("L141").Value = 70106.664
$ws.Range("N141").Value = -80466.664

$ws = $wb.Worksheets.Item("CUL")
# Row 68: Such a Butter Face
$ws.Range("H68").Value = 2112.951
$ws.Range("I68").Value = 839.8
$ws.Range("J68").Value = 2226.625
$ws.Range("K68").Value = 2519.4
$ws.Range("L68").Value = 6679.875
$ws.Range("M68").Value = -1708.4
$ws.Range("N68").Value = -8301.875

# Row 71: No Margarine of Error (L)
$ws.Range("H71").Value = 2112.951
$ws.Range("I71").Value = 839.8
$ws.Range("J71").Value = 2226.625
$ws.Range("K71").Value = 7558.2
$ws.Range("L71").Value = 20039.625
$ws.Range("M71").Value = -3502.2
$ws.Range("N71").Value = -28151.625

# Row 93: Loquacious
$ws.Range("H93").Value = 5187
$ws.Range("J93").Value = 5757.125
$ws.Range("L93").Value = 17271.375
$ws.Range("N93").Value = -21015.375

# Row 104: Fits to a Tea
$ws.Range("H104").Value = 6000
$ws.Range("J104").Value = 6000
$ws.Range("L104").Value = 18000
$ws.Range("N104").Value = -23242

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 15174092
$ws.Range("J131").Value = 29657.24
$ws.Range("L131").Value = 88971.72
$ws.Range("N131").Value = -99051.72

$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = 0

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 2648372.5
$ws.Range("J126").Value = 3004.6667
$ws.Range("L126").Value = 9014.000100000001
$ws.Range("N126").Value = -13954.0001

$ws = $wb.Worksheets.Item("LTW")
# Row 55: It's Not a Job, It's a Calling
$ws.Range("H55").Value = 180.83333
$ws.Range("I55").Value = 157
$ws.Range("J55").Value = 300
$ws.Range("K55").Value = 157
$ws.Range("L55").Value = 300
$ws.Range("M55").Value = 16
$ws.Range("N55").Value = -646

# Row 135: Dreams of Ja
$ws.Range("H135").Value = 32429
$ws.Range("J135").Value = 32429
$ws.Range("L135").Value = 32429
$ws.Range("N135").Value = -42569

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 3813.72
$ws.Range("I136").Value = 2266.125
$ws.Range("K136").Value = 6798.375
$ws.Range("M136").Value = -4248.375

$ws = $wb.Worksheets.Item("WVR")
# Row 62: Pride Up in Smoke
$ws.Range("H62").Value = 2287.5
$ws.Range("J62").Value = 2285
$ws.Range("L62").Value = 2285
$ws.Range("N62").Value = -3533

# Row 65: Desperate for Diversionaries (L)
$ws.Range("H65").Value = 2287.5
$ws.Range("J65").Value = 2285
$ws.Range("L65").Value = 11425
$ws.Range("N65").Value = -17665

# Row 122: Heavy Armoire
$ws.Range("H122").Value = 209959.22
$ws.Range("I122").Value = 235956.75
$ws.Range("K122").Value = 707870.25
$ws.Range("M122").Value = -705420.25
